# "version add product - tax"
# Adds a new "vat" column (M) to the AddProduct sheet with per-row tax
# values, and updates the sheet selection to reflect where the user
# ended up after adding the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the added "vat" column.
$ws.Range("M1").Value = "vat"

# Per-row vat values for the three data rows.
$ws.Range("M2").Value = 5
$ws.Range("M3").Value = 2
$ws.Range("M4").Value = 2

# Reflect the post-edit selection (next empty cell after the new column).
$ws.Range("N2").Select() | Out-Null
